# Auto-generated script applying cell value updates per the target diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2845.5454
$ws.Range("I43").Value = 5225.25
$ws.Range("J43").Value = 1485.7142
$ws.Range("K43").Value = 5225.25
$ws.Range("L43").Value = 1485.7142
$ws.Range("M43").Value = -5156.25
$ws.Range("N43").Value = -1623.7142
$ws.Range("H64").Value = 3983.3333
$ws.Range("I64").Value = 3700
$ws.Range("J64").Value = 5400
$ws.Range("K64").Value = 3700
$ws.Range("L64").Value = 5400
$ws.Range("M64").Value = -3452
$ws.Range("N64").Value = -5896
$ws.Range("H67").Value = 3983.3333
$ws.Range("I67").Value = 3700
$ws.Range("J67").Value = 5400
$ws.Range("K67").Value = 3700
$ws.Range("L67").Value = 5400
$ws.Range("M67").Value = -2842
$ws.Range("N67").Value = -7116
$ws.Range("H74").Value = 7279512.5
$ws.Range("I74").Value = 7279512.5
$ws.Range("K74").Value = 7279512.5
$ws.Range("M74").Value = -7278576.5
$ws.Range("H76").Value = 62502784
$ws.Range("I76").Value = 71431300
$ws.Range("J76").Value = 3166.6667
$ws.Range("K76").Value = 71431300
$ws.Range("L76").Value = 3166.6667
$ws.Range("M76").Value = -71430985
$ws.Range("N76").Value = -3796.6667
$ws.Range("H77").Value = 7279512.5
$ws.Range("I77").Value = 7279512.5
$ws.Range("K77").Value = 36397562.5
$ws.Range("M77").Value = -36392882.5
$ws.Range("H79").Value = 62502784
$ws.Range("I79").Value = 71431300
$ws.Range("J79").Value = 3166.6667
$ws.Range("K79").Value = 71431300
$ws.Range("L79").Value = 3166.6667
$ws.Range("M79").Value = -71430208
$ws.Range("N79").Value = -5350.6667
$ws.Range("H98").Value = 2161.9443
$ws.Range("I98").Value = 2387.2666
$ws.Range("K98").Value = 2387.2666
$ws.Range("M98").Value = -889.2665999999999
$ws.Range("H122").Value = 2161.9443
$ws.Range("I122").Value = 2387.2666
$ws.Range("K122").Value = 7161.7998
$ws.Range("M122").Value = -4711.7998
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 42224.88
$ws.Range("I2").Value = 57368.723
$ws.Range("K2").Value = 57368.723
$ws.Range("M2").Value = -57255.723
$ws.Range("H63").Value = 2247.4092
$ws.Range("I63").Value = 2156.647
$ws.Range("J63").Value = 2556
$ws.Range("K63").Value = 2156.647
$ws.Range("L63").Value = 2556
$ws.Range("M63").Value = -1470.647
$ws.Range("N63").Value = -3928
$ws.Range("H66").Value = 2247.4092
$ws.Range("I66").Value = 2156.647
$ws.Range("J66").Value = 2556
$ws.Range("K66").Value = 10783.235
$ws.Range("L66").Value = 12780
$ws.Range("M66").Value = -7351.235000000001
$ws.Range("N66").Value = -19644
$ws.Range("H88").Value = 3192
$ws.Range("I88").Value = 2740.3333
$ws.Range("J88").Value = 3417.8333
$ws.Range("K88").Value = 2740.3333
$ws.Range("L88").Value = 3417.8333
$ws.Range("M88").Value = -2334.3333
$ws.Range("N88").Value = -4229.8333
$ws.Range("H91").Value = 3192
$ws.Range("I91").Value = 2740.3333
$ws.Range("J91").Value = 3417.8333
$ws.Range("K91").Value = 2740.3333
$ws.Range("L91").Value = 3417.8333
$ws.Range("M91").Value = -1336.3333
$ws.Range("N91").Value = -6225.8333
$ws.Range("H116").Value = 42224.88
$ws.Range("I116").Value = 57368.723
$ws.Range("K116").Value = 57368.723
$ws.Range("M116").Value = -55074.723
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 42224.88
$ws.Range("I3").Value = 57368.723
$ws.Range("K3").Value = 57368.723
$ws.Range("M3").Value = -57254.723
$ws.Range("H86").Value = 34485176
$ws.Range("I86").Value = 58825696
$ws.Range("J86").Value = 2771.0833
$ws.Range("K86").Value = 58825696
$ws.Range("L86").Value = 2771.0833
$ws.Range("M86").Value = -58824573
$ws.Range("N86").Value = -5017.0833
$ws.Range("H89").Value = 34485176
$ws.Range("I89").Value = 58825696
$ws.Range("J89").Value = 2771.0833
$ws.Range("K89").Value = 294128480
$ws.Range("L89").Value = 13855.4165
$ws.Range("M89").Value = -294122864
$ws.Range("N89").Value = -25087.4165
$ws.Range("H105").Value = 3670.2046
$ws.Range("I105").Value = 3480.6667
$ws.Range("J105").Value = 3971.2354
$ws.Range("K105").Value = 3480.6667
$ws.Range("L105").Value = 3971.2354
$ws.Range("M105").Value = -1733.6667
$ws.Range("N105").Value = -7465.2354
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5500.7144
$ws.Range("I62").Value = 6502.5
$ws.Range("J62").Value = 5100
$ws.Range("K62").Value = 6502.5
$ws.Range("L62").Value = 5100
$ws.Range("M62").Value = -5878.5
$ws.Range("N62").Value = -6348
$ws.Range("H65").Value = 5500.7144
$ws.Range("I65").Value = 6502.5
$ws.Range("J65").Value = 5100
$ws.Range("K65").Value = 32512.5
$ws.Range("L65").Value = 25500
$ws.Range("M65").Value = -29392.5
$ws.Range("N65").Value = -31740
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1556.05
$ws.Range("I5").Value = 365.7857
$ws.Range("J5").Value = 4333.3335
$ws.Range("K5").Value = 1097.3571
$ws.Range("L5").Value = 13000.0005
$ws.Range("M5").Value = -985.3571000000002
$ws.Range("N5").Value = -13224.0005
$ws.Range("H80").Value = 7283.2144
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 7283.2144
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 21849.6432
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -23721.6432
$ws.Range("H83").Value = 7283.2144
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 7283.2144
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 65548.9296
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -74908.9296
$ws.Range("H128").Value = 989950
$ws.Range("I128").Value = 989950
$ws.Range("K128").Value = 2969850
$ws.Range("M128").Value = -2964870
$ws.Range("H132").Value = 1178.2354
$ws.Range("I132").Value = 855.3333
$ws.Range("J132").Value = 3600
$ws.Range("K132").Value = 7697.9997
$ws.Range("L132").Value = 32400
$ws.Range("M132").Value = -5167.9997
$ws.Range("N132").Value = -37460
$ws.Range("H133").Value = 10000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 10000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 30000
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -40120
$ws.Range("H134").Value = 4883.1055
$ws.Range("I134").Value = 3406.5833
$ws.Range("J134").Value = 7414.2856
$ws.Range("K134").Value = 10219.7499
$ws.Range("L134").Value = 22242.8568
$ws.Range("M134").Value = -5149.749899999999
$ws.Range("N134").Value = -32382.8568
$ws.Range("H135").Value = 1556.05
$ws.Range("I135").Value = 365.7857
$ws.Range("J135").Value = 4333.3335
$ws.Range("K135").Value = 3292.0713
$ws.Range("L135").Value = 39000.0015
$ws.Range("M135").Value = -757.0713000000001
$ws.Range("N135").Value = -44070.0015
$ws.Range("H136").Value = 55082.58
$ws.Range("I136").Value = 84480.75
$ws.Range("J136").Value = 4685.7144
$ws.Range("K136").Value = 253442.25
$ws.Range("L136").Value = 14057.1432
$ws.Range("M136").Value = -248342.25
$ws.Range("N136").Value = -24257.1432
$ws.Range("H137").Value = 4055.56
$ws.Range("I137").Value = 2999.2856
$ws.Range("J137").Value = 5399.909
$ws.Range("K137").Value = 8997.856800000001
$ws.Range("L137").Value = 16199.727
$ws.Range("M137").Value = -3897.856800000001
$ws.Range("N137").Value = -26399.727
$ws.Range("H138").Value = 1787.9231
$ws.Range("I138").Value = 1550.8334
$ws.Range("J138").Value = 4633
$ws.Range("K138").Value = 4652.5002
$ws.Range("L138").Value = 13899
$ws.Range("M138").Value = 487.4997999999996
$ws.Range("N138").Value = -24179
$ws.Range("H139").Value = 38131.035
$ws.Range("I139").Value = 39366.848
$ws.Range("K139").Value = 118100.544
$ws.Range("M139").Value = -112960.544
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4897
$ws.Range("I70").Value = 4835.5713
$ws.Range("J70").Value = 5004.5
$ws.Range("K70").Value = 4835.5713
$ws.Range("L70").Value = 5004.5
$ws.Range("M70").Value = -4565.5713
$ws.Range("N70").Value = -5544.5
$ws.Range("H73").Value = 4897
$ws.Range("I73").Value = 4835.5713
$ws.Range("J73").Value = 5004.5
$ws.Range("K73").Value = 4835.5713
$ws.Range("L73").Value = 5004.5
$ws.Range("M73").Value = -3899.5713
$ws.Range("N73").Value = -6876.5
$ws.Range("H132").Value = 2150.606
$ws.Range("I132").Value = 2114.1191
$ws.Range("J132").Value = 2214.4583
$ws.Range("K132").Value = 6342.3573
$ws.Range("L132").Value = 6643.374899999999
$ws.Range("M132").Value = -3812.3573
$ws.Range("N132").Value = -11703.3749
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 774.53845
$ws.Range("I132").Value = 664.7174
$ws.Range("J132").Value = 1616.5
$ws.Range("K132").Value = 1994.1522
$ws.Range("L132").Value = 4849.5
$ws.Range("M132").Value = 535.8478
$ws.Range("N132").Value = -9909.5
